# Weekly update: a new price-survey record (Asterix, "1a (cosecha)",
# 2022-01-28, Región del Maule) is inserted as row 80, pushing the
# existing rows 80-169 down to 81-170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 80; everything below shifts down.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(80, 1).Value = 11
$ws.Cells.Item(80, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(80, 3).Value = "Bíobío"
$ws.Cells.Item(80, 4).Value = 44589
$ws.Cells.Item(80, 5).Value = 8
$ws.Cells.Item(80, 6).Value = 100114001
$ws.Cells.Item(80, 7).Value = "Papa"
$ws.Cells.Item(80, 8).Value = "Asterix"
$ws.Cells.Item(80, 9).Value = "1a (cosecha)"
$ws.Cells.Item(80, 10).Value = 450
$ws.Cells.Item(80, 11).Value = 8000
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = 8556
$ws.Cells.Item(80, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Región del Maule"
$ws.Cells.Item(80, 16).Value = 342
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
